$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Durations_psplib")

$ws.Range("C1").Value = 157.631557226181

$ws.Range("A4").Value = 18813.769
$ws.Range("B4").Value = 18280
$ws.Range("F4").Value = 6412.029
$ws.Range("G4").Value = 6327
